# Automatic update of files.
# Applies the species-occurrence record shuffle + id bump described by the diff:
#   - Row 12 and Row 13 exchange their record-specific data (coords, accuracy,
#     start/end time, observer), row 13 gaining the time fields row 12 had and
#     row 12 losing them.
#   - Rows 15, 16 and 17 cyclically rotate their record-specific data
#     (15<-17, 16<-15, 17<-16).
#   - Column B ("Taxonsorteringsordning") bumps 79244 -> 79245 on every
#     touched row (9, 12, 13, 15, 16, 17, 19, 20).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 9 : only the sort-order id changes -------------------------------
$ws.Range("B9").Value2 = 79245

# ---- Rows 12 <-> 13 : full swap of the record-specific fields -------------
$ws.Range("A12").Value2 = 130979103
$ws.Range("Q12").Value2 = 570739
$ws.Range("R12").Value2 = 6736418
$ws.Range("S12").Value2 = 1
$ws.Range("Z12").ClearContents()
$ws.Range("AB12").ClearContents()
$ws.Range("AF12").ClearContents()
$ws.Range("AW12").Value2 = "Erik Danielsson"
$ws.Range("AX12").Value2 = "Erik Danielsson"
$ws.Range("B12").Value2 = 79245

$ws.Range("A13").Value2 = 130983072
$ws.Range("Q13").Value2 = 570809
$ws.Range("R13").Value2 = 6736404
$ws.Range("S13").Value2 = 10
$ws.Range("Z13").Value2 = "08:44"
$ws.Range("AB13").Value2 = "08:44"
$ws.Range("AF13").Value2 = ""
$ws.Range("AW13").Value2 = "Bo karlstens"
$ws.Range("AX13").Value2 = "Bo karlstens"
$ws.Range("B13").Value2 = 79245

# ---- Rows 15 -> 16 -> 17 -> 15 : cyclic rotation of record-specific fields -
$ws.Range("A15").Value2 = 130983074
$ws.Range("Q15").Value2 = 570764
$ws.Range("R15").Value2 = 6736425
$ws.Range("Z15").Value2 = "08:23"
$ws.Range("AB15").Value2 = "08:23"
$ws.Range("B15").Value2 = 79245

$ws.Range("A16").Value2 = 130983071
$ws.Range("P16").Value2 = "Flytjärnsmyran, Dlr"
$ws.Range("Q16").Value2 = 570817
$ws.Range("R16").Value2 = 6736417
$ws.Range("Z16").Value2 = "08:53"
$ws.Range("AB16").Value2 = "08:53"
$ws.Range("AF16").Value2 = ""
$ws.Range("AW16").Value2 = "Bo karlstens"
$ws.Range("AX16").Value2 = "Bo karlstens"
$ws.Range("B16").Value2 = 79245

$ws.Range("A17").Value2 = 130983619
$ws.Range("P17").Value2 = "Flytjärnsmyren, Dlr"
$ws.Range("Q17").Value2 = 570825
$ws.Range("R17").Value2 = 6736389
$ws.Range("Z17").Value2 = "08:54"
$ws.Range("AB17").Value2 = "08:54"
$ws.Range("AF17").ClearContents()
$ws.Range("AW17").Value2 = "Göran Ehn"
$ws.Range("AX17").Value2 = "Göran Ehn"
$ws.Range("B17").Value2 = 79245

# ---- Rows 19, 20 : only the sort-order id changes --------------------------
$ws.Range("B19").Value2 = 79245
$ws.Range("B20").Value2 = 79245
